$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.476312637329102
$ws.Range("B1").Value = 3.997977256774902
$ws.Range("C1").Value = 0.1788633912801743
$ws.Range("D1").Value = 0.1400669068098068
$ws.Range("E1").Value = 0.1275423914194107
